$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 737.6429000000001
$ws.Range("J41").Value = 278.4
$ws.Range("L41").Value = 278.4
$ws.Range("N41").Value = -1158.4

$ws.Range("H88").Value = 1669.25
$ws.Range("I88").Value = 3100
$ws.Range("J88").Value = 1192.3334
$ws.Range("K88").Value = 3100
$ws.Range("L88").Value = 1192.3334
$ws.Range("M88").Value = -2694
$ws.Range("N88").Value = -2004.3334

$ws.Range("H91").Value = 1669.25
$ws.Range("I91").Value = 3100
$ws.Range("J91").Value = 1192.3334
$ws.Range("K91").Value = 3100
$ws.Range("L91").Value = 1192.3334
$ws.Range("M91").Value = -1696
$ws.Range("N91").Value = -4000.3334

$ws.Range("H92").Value = 18320.643
$ws.Range("J92").Value = 460
$ws.Range("L92").Value = 460
$ws.Range("N92").Value = -2956

$ws.Range("H98").Value = 58995.25
$ws.Range("J98").Value = 6898.6665
$ws.Range("L98").Value = 6898.6665
$ws.Range("N98").Value = -9894.666499999999

$ws.Range("H99").Value = 286.0909
$ws.Range("I99").Value = 274.7
$ws.Range("J99").Value = 400
$ws.Range("K99").Value = 824.0999999999999
$ws.Range("L99").Value = 1200
$ws.Range("M99").Value = 673.9000000000001
$ws.Range("N99").Value = -4196

$ws.Range("H113").Value = 2487.4443
$ws.Range("I113").Value = 2065
$ws.Range("K113").Value = 2065
$ws.Range("M113").Value = 1189

$ws.Range("H122").Value = 58995.25
$ws.Range("J122").Value = 6898.6665
$ws.Range("L122").Value = 20695.9995
$ws.Range("N122").Value = -25595.9995

$ws.Range("H132").Value = 1442
$ws.Range("I132").Value = 1574.9578
$ws.Range("K132").Value = 4724.873399999999
$ws.Range("M132").Value = -2194.873399999999

$ws.Range("H135").Value = 36794.08
$ws.Range("J135").Value = 674.5
$ws.Range("L135").Value = 6070.5
$ws.Range("N135").Value = -11140.5

$ws.Range("H138").Value = 1966.0312
$ws.Range("I138").Value = 1212.3695
$ws.Range("K138").Value = 3637.1085
$ws.Range("M138").Value = 1502.8915

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1863.6154
$ws.Range("I2").Value = 2043.9
$ws.Range("J2").Value = 1262.6666
$ws.Range("K2").Value = 2043.9
$ws.Range("L2").Value = 1262.6666
$ws.Range("M2").Value = -1930.9
$ws.Range("N2").Value = -1488.6666

$ws.Range("H97").Value = 636
$ws.Range("I97").Value = 387.93103
$ws.Range("J97").Value = 1355.4
$ws.Range("K97").Value = 387.93103
$ws.Range("L97").Value = 1355.4
$ws.Range("M97").Value = 108.06897
$ws.Range("N97").Value = -2347.4

$ws.Range("H116").Value = 1863.6154
$ws.Range("I116").Value = 2043.9
$ws.Range("J116").Value = 1262.6666
$ws.Range("K116").Value = 2043.9
$ws.Range("L116").Value = 1262.6666
$ws.Range("M116").Value = 250.0999999999999
$ws.Range("N116").Value = -5850.6666

$ws.Range("H122").Value = 72881.69500000001
$ws.Range("I122").Value = 3562.1428
$ws.Range("K122").Value = 10686.4284
$ws.Range("M122").Value = -8236.428400000001

$ws.Range("H132").Value = 8094.1177
$ws.Range("I132").Value = 13035.577
$ws.Range("J132").Value = 2955
$ws.Range("K132").Value = 39106.731
$ws.Range("L132").Value = 8865
$ws.Range("M132").Value = -36576.731
$ws.Range("N132").Value = -13925

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1863.6154
$ws.Range("I3").Value = 2043.9
$ws.Range("J3").Value = 1262.6666
$ws.Range("K3").Value = 2043.9
$ws.Range("L3").Value = 1262.6666
$ws.Range("M3").Value = -1929.9
$ws.Range("N3").Value = -1490.6666

$ws.Range("H103").Value = 15043.733
$ws.Range("J103").Value = 15043.733
$ws.Range("L103").Value = 15043.733
$ws.Range("N103").Value = -17387.733

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2391.2222
$ws.Range("I31").Value = 1977.1025
$ws.Range("J31").Value = 2775.762
$ws.Range("K31").Value = 1977.1025
$ws.Range("L31").Value = 2775.762
$ws.Range("M31").Value = -1682.1025
$ws.Range("N31").Value = -3365.762

$ws.Range("H34").Value = 2391.2222
$ws.Range("I34").Value = 1977.1025
$ws.Range("J34").Value = 2775.762
$ws.Range("K34").Value = 1977.1025
$ws.Range("L34").Value = 2775.762
$ws.Range("M34").Value = -1775.1025
$ws.Range("N34").Value = -3179.762

$ws.Range("H134").Value = 2248.6775
$ws.Range("I134").Value = 2344.3264
$ws.Range("K134").Value = 7032.9792
$ws.Range("M134").Value = -4497.9792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 721.26086
$ws.Range("I5").Value = 668.3684
$ws.Range("K5").Value = 2005.1052
$ws.Range("M5").Value = -1893.1052

$ws.Range("H116").Value = 13670.5
$ws.Range("I116").Value = 14505.75
$ws.Range("K116").Value = 43517.25
$ws.Range("M116").Value = -40075.25

$ws.Range("H119").Value = 7771.2
$ws.Range("I119").Value = 7771.2
$ws.Range("K119").Value = 23313.6
$ws.Range("M119").Value = -18475.6

$ws.Range("H135").Value = 721.26086
$ws.Range("I135").Value = 668.3684
$ws.Range("K135").Value = 6015.3156
$ws.Range("M135").Value = -3480.3156

$ws.Range("H138").Value = 8471.75
$ws.Range("I138").Value = 6944.75
$ws.Range("K138").Value = 20834.25
$ws.Range("M138").Value = -15694.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 39990
$ws.Range("J86").Value = 39990
$ws.Range("L86").Value = 39990
$ws.Range("N86").Value = -42362

$ws.Range("H89").Value = 39990
$ws.Range("J89").Value = 39990
$ws.Range("L89").Value = 119970
$ws.Range("N89").Value = -131826

$ws.Range("H97").Value = 15166508
$ws.Range("I97").Value = 19249328
$ws.Range("K97").Value = 19249328
$ws.Range("M97").Value = -19248832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 65966.664
$ws.Range("J6").Value = 65966.664
$ws.Range("L6").Value = 65966.664
$ws.Range("N6").Value = -66190.664

$ws.Range("H7").Value = 14831.333
$ws.Range("I7").Value = 16125.875
$ws.Range("K7").Value = 16125.875
$ws.Range("M7").Value = -16013.875

$ws.Range("H93").Value = 982
$ws.Range("I93").Value = 982
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 982
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 266
$ws.Range("N93").ClearContents()

$ws.Range("H96").Value = 31999
$ws.Range("J96").Value = 31999
$ws.Range("L96").Value = 31999
$ws.Range("N96").Value = -37491

$ws.Range("H122").Value = 5979.4287
$ws.Range("I122").Value = 3400.9614
$ws.Range("K122").Value = 10202.8842
$ws.Range("M122").Value = -7752.8842

$ws.Range("H126").Value = 14831.333
$ws.Range("I126").Value = 16125.875
$ws.Range("K126").Value = 48377.625
$ws.Range("M126").Value = -45907.625

$ws.Range("H132").Value = 2902.8125
$ws.Range("I132").Value = 2509.6667
$ws.Range("J132").Value = 4082.25
$ws.Range("K132").Value = 7529.000100000001
$ws.Range("L132").Value = 12246.75
$ws.Range("M132").Value = -4999.000100000001
$ws.Range("N132").Value = -17306.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 47619744
$ws.Range("I100").Value = 58824250
$ws.Range("J100").Value = 605.25
$ws.Range("K100").Value = 117648500
$ws.Range("L100").Value = 1210.5
$ws.Range("M100").Value = -117647959
$ws.Range("N100").Value = -2292.5

$ws.Range("H113").Value = 855.03705
$ws.Range("I113").Value = 917.86957
$ws.Range("J113").Value = 493.75
$ws.Range("K113").Value = 2753.60871
$ws.Range("L113").Value = 1481.25
$ws.Range("M113").Value = -583.60871
$ws.Range("N113").Value = -5821.25

$ws.Range("H132").Value = 2504.275
$ws.Range("I132").Value = 2540.7334
$ws.Range("J132").Value = 2394.9
$ws.Range("K132").Value = 7622.2002
$ws.Range("L132").Value = 7184.700000000001
$ws.Range("M132").Value = -5092.2002
$ws.Range("N132").Value = -12244.7
